# DeveloperGuide: Updated diagram & description for Logic, Model, Storage &
# Versioned Tasketch to suit our project.
#
# Renames the "AddressBook" domain used throughout the Undo/Redo sequence
# diagram to "TaskBook" (the app being documented changed from an address
# book to a task manager), plus a couple of small formatting tweaks that
# came along with the same PowerPoint edit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. ":AddressBookParser" -> ":TaskBookParser" (shape id=16) ----------
# Text is split across two paragraphs (":Address" / "BookParser"); only the
# "Address" part of the first paragraph changes.
$sh = $s.Shapes.Item(6)
$tr = $sh.TextFrame.TextRange
$idx = $tr.Text.IndexOf("Address")
$tr.Characters($idx + 1, 7).Text = "Task"

# --- 2. "undoAddressBook()" -> "undoTaskBook()" (shape id=79) ------------
$sh = $s.Shapes.Item(19)
$tr = $sh.TextFrame.TextRange
$idx = $tr.Text.IndexOf("Address")
$tr.Characters($idx + 1, 7).Text = "Task"

# --- 3. ":VersionedAddressBook" -> ":VersionedTaskBook" (shape id=84) ----
$sh = $s.Shapes.Item(23)
$tr = $sh.TextFrame.TextRange
$idx = $tr.Text.IndexOf("Address")
$tr.Characters($idx + 1, 7).Text = "Task"

# --- 4. TextBox 3 (shape id=4, the little "X"): wrap="none" -> "square" --
$sh = $s.Shapes.Item(30)
$sh.TextFrame.WordWrap = -1

# --- 5. "resetData(ReadOnlyAddressBook)" -> "resetData(ReadOnlyTaskBook)",
#        plus explicit left alignment on that paragraph (shape id=88) -----
$sh = $s.Shapes.Item(35)
$tr = $sh.TextFrame.TextRange
$tr.ParagraphFormat.Alignment = 1
$idx = $tr.Text.IndexOf("Address")
$tr.Characters($idx + 1, 7).Text = "Task"

# --- 6. Refresh the cached "today" date field text wherever it is still
#        reachable through the object model (slide master + every slide
#        layout). 7/6/2018 -> 09-Apr-19 ------------------------------------
$newDate = "09-Apr-19"
$oldDate = "7/6/2018"

function Update-DateShape($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

Update-DateShape($p.SlideMaster)

$layouts = $p.SlideMaster.CustomLayouts
for ($l = 1; $l -le $layouts.Count; $l++) {
    Update-DateShape($layouts.Item($l))
}
